# MALS-1042 Updating template to match view changes
#
# The envelope's "Client" address line merges the registrant's name as
# two adjacent fields: {d.RegistrantFirst}{d.Last}. The view now builds
# the name the other way around, so the template text needs to read
# {d.Registrant}{d.LastFirst} instead -- i.e. "RegistrantFirst" becomes
# "Registrant" and "Last" becomes "LastFirst".
$d = $word.ActiveDocument

$d.Content.Find.Execute("RegistrantFirstLast", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "RegistrantLastFirst", 2) | Out-Null
